# Santhosh L NestJS Time Sheet Attribute.xlsx
# Add files via upload
#
# Row 19 (03-Feb-2025 / 45603) previously only had the date in column A.
# This edit fills in the Task Name / Status columns for that row with the
# same values already used on the row above (row 18): "Krion 6D Help
# Document Modification" / "Onprocess". Also moves the sheet's active
# selection to B22 (scrolled back to the top of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = "Krion 6D Help Document Modification"
$ws.Range("C19").Value = "Onprocess"

$ws.Range("B22").Select()
